$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 148, shifting the existing rows (old 148-159) down to 150-161.
$ws.Rows("148:149").Insert()

# New row 148: Black Amber, Primera
$ws.Cells.Item(148, 1).Value = 10
$ws.Cells.Item(148, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(148, 3).Value = "La Araucanía"
$ws.Cells.Item(148, 4).Value = 44578
$ws.Cells.Item(148, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(148, 5).Value = 9
$ws.Cells.Item(148, 6).Value = "Fruta"
$ws.Cells.Item(148, 7).Value = 100103
$ws.Cells.Item(148, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(148, 9).Value = 100103002
$ws.Cells.Item(148, 10).Value = "Ciruela"
$ws.Cells.Item(148, 11).Value = "Black Amber"
$ws.Cells.Item(148, 12).Value = "Primera"
$ws.Cells.Item(148, 13).Value = 95
$ws.Cells.Item(148, 14).Value = 18000
$ws.Cells.Item(148, 15).Value = 18000
$ws.Cells.Item(148, 16).Value = 18000
$ws.Cells.Item(148, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(148, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(148, 19).Value = 1000
$ws.Cells.Item(148, 20).Value = 18

# New row 149: Crimsom fall, Primera
$ws.Cells.Item(149, 1).Value = 10
$ws.Cells.Item(149, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(149, 3).Value = "La Araucanía"
$ws.Cells.Item(149, 4).Value = 44578
$ws.Cells.Item(149, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(149, 5).Value = 9
$ws.Cells.Item(149, 6).Value = "Fruta"
$ws.Cells.Item(149, 7).Value = 100103
$ws.Cells.Item(149, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(149, 9).Value = 100103002
$ws.Cells.Item(149, 10).Value = "Ciruela"
$ws.Cells.Item(149, 11).Value = "Crimsom fall"
$ws.Cells.Item(149, 12).Value = "Primera"
$ws.Cells.Item(149, 13).Value = 110
$ws.Cells.Item(149, 14).Value = 12000
$ws.Cells.Item(149, 15).Value = 12000
$ws.Cells.Item(149, 16).Value = 12000
$ws.Cells.Item(149, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(149, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(149, 19).Value = 667
$ws.Cells.Item(149, 20).Value = 18
